$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "-"

# Row 3
$ws.Range("E3").Value = "[-, 'MCT-3A-Eletropneumática', -, -]"

# Row 7
$ws.Range("E7").Value = "[-, 'MCT-3A-Eletropneumática', -, -]"

# Row 8
$ws.Range("D8").Value = "-"

# Row 10
$ws.Range("E10").Value = "-"

# Row 11
$ws.Range("F11").Value = "[-, -, 'MEC-3A-C.pneumática', -]"

# Row 12
$ws.Range("F12").Value = "[-, -, 'MEC-3A-C.pneumática', -]"

# Row 14
$ws.Range("F14").Value = "[-, -, 'MEC-3A-C.pneumática', -]"

# Row 15
$ws.Range("F15").Value = "[-, -, 'MEC-3A-C.pneumática', -]"

# Row 16
$ws.Range("E16").Value = "-"

# Row 18
$ws.Range("B18").Value = "['MEC-1NB-Desenho tecnico mecanico', -]"
$ws.Range("D18").Value = "['MEC-1NA-Desenho tecnico mecanico – T1', -]"

# Row 19
$ws.Range("B19").Value = "['MEC-1NB-Desenho tecnico mecanico', -]"
$ws.Range("C19").Value = "-"
$ws.Range("E19").Value = "MEC-1NB-M.T.F."

# Row 20
$ws.Range("B20").Value = "['MEC-1NB-Desenho tecnico mecanico', -]"
$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "['MEC-1NA-Desenho tecnico mecanico – T1', -]"
$ws.Range("E20").Value = "MEC-1NB-M.T.F."

# Row 21
$ws.Range("B21").Value = "-"
$ws.Range("D21").Value = "['MEC-1NA-Desenho tecnico mecanico – T1', -]"
$ws.Range("E21").Value = "-"
